# Insert a new data row at row 22 (pushing the existing rows 22-58 down to 23-59)
# and populate it with the new weekly price record, per the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 22.. down by one, inheriting formatting from row 22.
$ws.Rows(22).Insert()

# Populate the newly inserted row 22 with the new record.
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44519
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 100112052
$ws.Range("G22").Value = "Albahaca"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 3500
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = 3750
$ws.Range("N22").Value = "$/paquete"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 3750
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
